$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-14 18:16:05", 0.0004),
    @("2023-12-14 18:16:34", 0.0016),
    @("2023-12-14 18:17:02", 0.002),
    @("2023-12-14 18:17:08", 0.0004)
)

$startRow = 331
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
